# Repull data, push all data, mean calculation
# Update the "dSF" (delta score final) column (F) for rows where the
# final data differs from the initially pulled data (column E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    4  = 2
    6  = -1
    15 = 3
    18 = -4
    22 = -2
    24 = -1
    32 = -5
    34 = 0
    35 = 2
    49 = 0
    54 = -4
    55 = 0
    56 = -2
    58 = 0
    60 = 4
    71 = 6
    73 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
